$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily log. It belongs right before the
# existing row 169, so insert a fresh row there (shifting rows 169:219 down
# to 170:220) and populate it with the new reading.
$ws.Rows(169).Insert()

$ws.Range('A169').Value2 = 11
$ws.Range('B169').Value2 = 'Vega Monumental Concepción'
$ws.Range('C169').Value2 = 'Bíobío'
$ws.Range('D169').Value2 = 44855
$ws.Range('E169').Value2 = 8
$ws.Range('F169').Value2 = 100112003
$ws.Range('G169').Value2 = 'Ajo'
$ws.Range('H169').Value2 = 'Chino'
$ws.Range('I169').Value2 = 'Primera'
$ws.Range('J169').Value2 = 400
$ws.Range('K169').Value2 = 14000
$ws.Range('L169').Value2 = 15000
$ws.Range('M169').Value2 = 14500
$ws.Range('N169').Value2 = '$/caja 10 kilos'
$ws.Range('O169').Value2 = 'China'
$ws.Range('P169').Value2 = 1450
$ws.Range('Q169').Value2 = 10
$ws.Range('R169').Value2 = 'Hortaliza'
